# regen save_data to use K instead of Strike#, regen std/mean, calc and write s_vals
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newK = @{
    2  = 1
    3  = 5
    4  = 3
    5  = 1
    6  = 7
    7  = 1
    8  = 4
    9  = 3
    10 = 6
    11 = 4
    12 = 3
    13 = 6
    14 = 8
    15 = 3
    16 = 4
    17 = 4
    18 = 5
    19 = 4
    20 = 4
}

foreach ($row in $newK.Keys) {
    $ws.Range("G$row").Value = $newK[$row]
}
